# Insert 3 new weekly price rows (date 44466) right before the current
# row 213, pushing all subsequent rows down by 3. Then populate the three
# new rows with the "Especial" / "Primera" / "Segunda" quality records for
# Packham's Triumph pears for that date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 213 (shifts old rows 213.. down to 216..)
$ws.Range("A213:A215").EntireRow.Insert()

# Row 213 - Packham's Triumph / Especial
$ws.Cells.Item(213, 1).Value = 8
$ws.Cells.Item(213, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(213, 3).Value = "Coquimbo"
$ws.Cells.Item(213, 4).Value = 44466
$ws.Cells.Item(213, 5).Value = 4
$ws.Cells.Item(213, 6).Value = "Fruta"
$ws.Cells.Item(213, 7).Value = 100104
$ws.Cells.Item(213, 8).Value = "Frutos de pepita"
$ws.Cells.Item(213, 9).Value = 100104005
$ws.Cells.Item(213, 10).Value = "Pera"
$ws.Cells.Item(213, 11).Value = "Packham's Triumph"
$ws.Cells.Item(213, 12).Value = "Especial"
$ws.Cells.Item(213, 13).Value = 16
$ws.Cells.Item(213, 14).Value = 285000
$ws.Cells.Item(213, 15).Value = 290000
$ws.Cells.Item(213, 16).Value = 287500
$ws.Cells.Item(213, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(213, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(213, 19).Value = 639
$ws.Cells.Item(213, 20).Value = 450

# Row 214 - Packham's Triumph / Primera
$ws.Cells.Item(214, 1).Value = 8
$ws.Cells.Item(214, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(214, 3).Value = "Coquimbo"
$ws.Cells.Item(214, 4).Value = 44466
$ws.Cells.Item(214, 5).Value = 4
$ws.Cells.Item(214, 6).Value = "Fruta"
$ws.Cells.Item(214, 7).Value = 100104
$ws.Cells.Item(214, 8).Value = "Frutos de pepita"
$ws.Cells.Item(214, 9).Value = 100104005
$ws.Cells.Item(214, 10).Value = "Pera"
$ws.Cells.Item(214, 11).Value = "Packham's Triumph"
$ws.Cells.Item(214, 12).Value = "Primera"
$ws.Cells.Item(214, 13).Value = 16
$ws.Cells.Item(214, 14).Value = 255000
$ws.Cells.Item(214, 15).Value = 260000
$ws.Cells.Item(214, 16).Value = 257500
$ws.Cells.Item(214, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(214, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(214, 19).Value = 572
$ws.Cells.Item(214, 20).Value = 450

# Row 215 - Packham's Triumph / Segunda
$ws.Cells.Item(215, 1).Value = 8
$ws.Cells.Item(215, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(215, 3).Value = "Coquimbo"
$ws.Cells.Item(215, 4).Value = 44466
$ws.Cells.Item(215, 5).Value = 4
$ws.Cells.Item(215, 6).Value = "Fruta"
$ws.Cells.Item(215, 7).Value = 100104
$ws.Cells.Item(215, 8).Value = "Frutos de pepita"
$ws.Cells.Item(215, 9).Value = 100104005
$ws.Cells.Item(215, 10).Value = "Pera"
$ws.Cells.Item(215, 11).Value = "Packham's Triumph"
$ws.Cells.Item(215, 12).Value = "Segunda"
$ws.Cells.Item(215, 13).Value = 20
$ws.Cells.Item(215, 14).Value = 235000
$ws.Cells.Item(215, 15).Value = 240000
$ws.Cells.Item(215, 16).Value = 237500
$ws.Cells.Item(215, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(215, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(215, 19).Value = 528
$ws.Cells.Item(215, 20).Value = 450
